$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "67.847.42"
Set-TextCell "E2" "  -4.59%  "
Set-TextCell "D3" "3.259.26"
Set-TextCell "E3" "  -8.12%  "
Set-TextCell "D5" "592.76"
Set-TextCell "E5" "  -4.22%  "
Set-TextCell "D6" "154.62"
Set-TextCell "E6" "  -11.27%  "
Set-TextCell "D7" "0.998"
Set-TextCell "E7" "  -0.13%  "
Set-TextCell "D8" "3.248.73"
Set-TextCell "E8" "  -8.26%  "
Set-TextCell "D9" "0.548"
Set-TextCell "E9" "  -10.55%  "
Set-TextCell "D10" "0.176"
Set-TextCell "E10" "  -11.71%  "
Set-TextCell "D11" "6.87"
Set-TextCell "E11" "  -4.99%  "
Set-TextCell "D12" "0.509"
Set-TextCell "E12" "  -13.87%  "
Set-TextCell "D13" "38.90"
Set-TextCell "E13" "  -16.89%  "
Set-TextCell "D14" "0.0000247"
Set-TextCell "E14" "  -10.99%  "
Set-TextCell "D15" "3.766.97"
Set-TextCell "E15" "  -8.40%  "
Set-TextCell "D16" "67.640.05"
Set-TextCell "E16" "  -4.83%  "
Set-TextCell "B17" "BitcoinCash"
Set-TextCell "C17" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D17" "546.27"
Set-TextCell "E17" "  -11.14%  "
Set-TextCell "B18" "WrappedEther"
Set-TextCell "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D18" "3.242.12"
Set-TextCell "E18" "  -8.52%  "
Set-TextCell "B19" "Polkadot"
Set-TextCell "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D19" "7.28"
Set-TextCell "E19" "  -14.06%  "
Set-TextCell "B20" "TRON"
Set-TextCell "C20" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D20" "0.115"
Set-TextCell "E20" "  -5.65%  "
Set-TextCell "D21" "15.30"
Set-TextCell "E21" "  -14.28%  "
Set-TextCell "D22" "0.770"
Set-TextCell "E22" "  -13.66%  "
Set-TextCell "D23" "7.88"
Set-TextCell "E23" "  -13.46%  "
Set-TextCell "D24" "86.25"
Set-TextCell "E24" "  -12.21%  "
Set-TextCell "D25" "13.69"
Set-TextCell "E25" "  -13.30%  "
Set-TextCell "E26" "  +0.03%  "
Set-TextCell "D27" "3.20"
Set-TextCell "E27" "  -15.83%  "
Set-TextCell "D28" "8.24"
Set-TextCell "E28" "  -10.47%  "
Set-TextCell "D29" "29.74"
Set-TextCell "E29" "  -12.43%  "
Set-TextCell "E30" "  -17.67%  "
Set-TextCell "D31" "2.73"
Set-TextCell "E31" "  -10.57%  "
Set-TextCell "E32" "  -11.28%  "
Set-TextCell "D33" "552.21"
Set-TextCell "E33" "  -10.30%  "
Set-TextCell "D34" "6.66"
Set-TextCell "E34" "  -18.82%  "
Set-TextCell "D35" "5.83"
Set-TextCell "E35" "  -15.39%  "
Set-TextCell "D36" "0.999"
Set-TextCell "E36" "  -0.11%  "
Set-TextCell "D37" "0.0449"
Set-TextCell "E37" "  -6.19%  "
Set-TextCell "D38" "54.12"
Set-TextCell "E38" "  -5.10%  "
Set-TextCell "D39" "0.0856"
Set-TextCell "E39" "  -15.12%  "
Set-TextCell "D40" "9.29"
Set-TextCell "E40" "  -14.47%  "
Set-TextCell "E41" "  -12.81%  "
Set-TextCell "D42" "2.954.01"
Set-TextCell "E42" "  -12.66%  "
Set-TextCell "D43" "2.67"
Set-TextCell "E43" "  -24.05%  "
Set-TextCell "D44" "0.0₃0598"
Set-TextCell "E44" "  -19.65%  "
Set-TextCell "D45" "0.264"
Set-TextCell "E45" "  -16.50%  "
Set-TextCell "D46" "2.41"
Set-TextCell "E46" "  -19.32%  "
Set-TextCell "D47" "26.40"
Set-TextCell "E47" "  -18.41%  "
Set-TextCell "B48" "Fetch.AI"
Set-TextCell "C48" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D48" "2.15"
Set-TextCell "E48" "  -16.51%  "
Set-TextCell "B49" "USDe"
Set-TextCell "C49" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D49" "1.00"
Set-TextCell "E49" "  -0.03%  "
Set-TextCell "E50" "  -12.75%  "
Set-TextCell "D51" "125.56"
Set-TextCell "E51" "  -6.01%  "
